## faux registry log modal, some tweaks
## Adds new "Registry" related localization strings to the language sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row right after the existing "search_title_files" row ---
# Before: row 19 = search_button / SEARCH
# After insertion, row 19 becomes the new registry search-title row, and
# everything that was there shifts down by one.
$ws.Rows(19).Insert()

# --- Append the new registry-inspector strings at the bottom of the table ---
# After the insert above, the former last row (39) is now row 40, so the
# new rows start at 41. Cells are written in the order the new shared
# strings were originally authored so the string table layout matches.
$ws.Range("A43").Value = "registry_log_title"
$ws.Range("B43").Value = "REGISTRY: Computer/HKEY_CURRENT_USER/Software/Macrohard/Doors/CurrentVersion/Run"

$ws.Range("A44").Value = "name"
$ws.Range("B44").Value = "Name"

$ws.Range("A45").Value = "type"
$ws.Range("B45").Value = "Type"

$ws.Range("A46").Value = "data"
$ws.Range("B46").Value = "Data"

$ws.Range("A42").Value = "registry_inspector"
$ws.Range("B42").Value = "Registry Inspector"

$ws.Range("A19").Value = "search_title_registry"
$ws.Range("B19").Value = "Search In Registry: {0}"

$ws.Range("A41").Value = "registry"
$ws.Range("B41").Value = "Registry"

# --- Update the view state to match where the author left off editing ---
$ws.Activate() | Out-Null
try {
  $excel.ActiveWindow.ScrollRow = 24
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("B41").Select() | Out-Null
